$d = $word.ActiveDocument

# Word Find.Execute positional parameters:
# (FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#  MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
# Wrap = 1 (wdFindContinue), Replace = 2 (wdReplaceAll)

# 1) "escolher sua temporada de" -> "escolher seu protagonista de"
#    (paragraph about the login screen / "The Walking Dead" season selection)
$range1 = $d.Content
$found1 = $range1.Find.Execute(
    "escolher sua temporada de", $false, $false, $false, $false, $false,
    $true, 1, $false, "escolher seu protagonista de", 2)
if (-not $found1) {
    throw "Could not find text 'escolher sua temporada de'"
}

# 2) "Walking Dead preferida" -> "Walking Dead game preferida"
#    (insert the word "game" before "preferida")
$range2 = $d.Content
$found2 = $range2.Find.Execute(
    "Walking Dead preferida", $false, $false, $false, $false, $false,
    $true, 1, $false, "Walking Dead game preferida", 2)
if (-not $found2) {
    throw "Could not find text 'Walking Dead preferida'"
}

# 3) "cards das temporadas. " -> "cards dos personagens "
#    (paragraph about the signup screen / cards visualization)
$range3 = $d.Content
$found3 = $range3.Find.Execute(
    "cards das temporadas. ", $false, $false, $false, $false, $false,
    $true, 1, $false, "cards dos personagens ", 2)
if (-not $found3) {
    throw "Could not find text 'cards das temporadas. '"
}

Write-Host "Replacements applied: $found1, $found2, $found3"
